$d = $word.ActiveDocument

# The template currently reads "Claim number: <<caseNumber>>" in the
# top-right address block of the first page. Update the wording so it
# reads "Case number: <<caseNumber>>" (the merge field itself, i.e. the
# "caseNumber" placeholder, is left untouched).
#
# We scope the search to the exact phrase "Claim number: <<c" (rather
# than just "Claim") so that unrelated text elsewhere in the document
# (e.g. "Claimant", "Claimant 1", "Claimant 2") is never touched.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$replaced = $d.Content.Find.Execute(
    "Claim number: <<c",  # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    "Case number: <<c",    # ReplaceWith
    2                       # Replace (wdReplaceOne)
)

Write-Host "Replaced claim->case heading: $replaced"
